$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New price (column D) values per row; rows not listed here keep their existing price.
$prices = @{
    2 = "272.10"
    3 = "23.24"
    4 = "6.373"
    5 = "0.06284"
    6 = "3.647"
    7 = "6.775"
    8 = "1.386"
    9 = "0.8363"
    10 = "0.1629"
    12 = "0.03495"
    13 = "0.03172"
    14 = "0.09308"
    15 = "3.943"
    16 = "0.001698"
    17 = "0.04857"
    18 = "0.006297"
    20 = "0.001088"
    21 = "0.0001499"
    22 = "3.733"
    23 = "2.311"
    24 = "0.01391"
    26 = "0.1262"
    27 = "0.0003735"
    40 = "0.04691"
    41 = "0.006890"
    42 = "0.1177"
    43 = "0.003456"
    45 = "0.00006265"
    46 = "0.00000000750"
    47 = "0.7968"
    48 = "0.1029"
    49 = "0.00002099"
    50 = "0.01240"
}

# Every data row (2-51) gets its Hora (column G) bumped from 15 to 16.
for ($row = 2; $row -le 51; $row++) {
    if ($prices.ContainsKey($row)) {
        $dCell = $ws.Cells.Item($row, 4)
        $dCell.NumberFormat = "@"
        $dCell.Value = $prices[$row]
    }
    $gCell = $ws.Cells.Item($row, 7)
    $gCell.NumberFormat = "@"
    $gCell.Value = "16"
}
